$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: "Latent variable models" / "latent-var-models"
$ws.Range("F14").Value = "Latent variable models"
$ws.Range("G14").Value = "latent-var-models"
$ws.Range("H14").Formula = '=CONCATENATE(G14,"-",SUBSTITUTE(A14,".",""))'

# Row 15: "Wekalink" / "wekalink"
$ws.Range("F15").Value = "Wekalink"
$ws.Range("G15").Value = "wekalink"
$ws.Range("H15").Formula = '=CONCATENATE(G15,"-",SUBSTITUTE(A15,".",""))'

# Update the active selection to K20
$ws.Range("K20").Select()
